$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.943.70'
$ws.Range("E2").Value = '  +0.74%  '

$ws.Range("D3").Value = '2.233.01'
$ws.Range("E3").Value = '  -3.27%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '''292.31'
$ws.Range("E5").Value = '  -4.60%  '

$ws.Range("D6").Value = '''86.17'
$ws.Range("E6").Value = '  +5.60%  '

$ws.Range("D7").Value = '''0.512'
$ws.Range("E7").Value = '  -0.87%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").Value = '''0.468'
$ws.Range("E9").Value = '  -0.42%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '''0.0798'
$ws.Range("E10").Value = '  +1.83%  '

$ws.Range("B11").Value = 'Avalanche'
$ws.Range("C11").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D11").Value = '''30.59'
$ws.Range("E11").Value = '  +6.77%  '

$ws.Range("D12").Value = '''47.03'
$ws.Range("E12").Value = '  -10.07%  '

$ws.Range("E13").Value = '  -1.18%  '

$ws.Range("E14").Value = '  +2.85%  '

$ws.Range("D15").Value = '2.580.07'
$ws.Range("E15").Value = '  -3.40%  '

$ws.Range("D16").Value = '''14.12'
$ws.Range("E16").Value = '  -1.87%  '

$ws.Range("D17").Value = '2.239.38'
$ws.Range("E17").Value = '  -3.86%  '

$ws.Range("D18").Value = '''0.725'
$ws.Range("E18").Value = '  -1.59%  '

$ws.Range("D19").Value = '39.863.61'
$ws.Range("E19").Value = '  +0.70%  '

$ws.Range("D20").Value = '0.0₃0893'
$ws.Range("E20").Value = '  +1.41%  '

$ws.Range("D21").Value = '''5.78'
$ws.Range("E21").Value = '  -2.33%  '

$ws.Range("D22").Value = '''10.64'
$ws.Range("E22").Value = '  +4.05%  '

$ws.Range("D23").Value = '''65.38'
$ws.Range("E23").Value = '  -2.87%  '

$ws.Range("D24").Value = '''235.04'
$ws.Range("E24").Value = '  +1.43%  '

$ws.Range("E25").Value = '  +0.09%  '

$ws.Range("D26").Value = '''2.42'
$ws.Range("E26").Value = '  -2.02%  '

$ws.Range("D27").Value = '''1.83'
$ws.Range("E27").Value = '  +3.87%  '

$ws.Range("D28").Value = '''22.91'
$ws.Range("E28").Value = '  +0.67%  '

$ws.Range("E29").Value = '  +1.46%  '

$ws.Range("E30").Value = '  +2.56%  '

$ws.Range("D31").Value = '''33.84'
$ws.Range("E31").Value = '  +4.40%  '

$ws.Range("D32").Value = '''154.36'
$ws.Range("E32").Value = '  +2.55%  '

$ws.Range("E33").Value = '  -0.22%  '

$ws.Range("E34").Value = '  -1.67%  '

$ws.Range("D35").Value = '''0.0709'
$ws.Range("E35").Value = '  +1.48%  '

$ws.Range("E36").Value = '  -2.49%  '

$ws.Range("D37").Value = '''16.46'
$ws.Range("E37").Value = '  +9.02%  '

$ws.Range("E38").Value = '  -0.16%  '

$ws.Range("E39").Value = '  +4.03%  '

$ws.Range("E40").Value = '  +0.09%  '

$ws.Range("D41").Value = '''1.66'
$ws.Range("E41").Value = '  +2.14%  '

$ws.Range("E42").Value = '  +3.17%  '

$ws.Range("D43").Value = '1.959.05'
$ws.Range("E43").Value = '  +0.03%  '

$ws.Range("E44").Value = '  -3.13%  '

$ws.Range("D45").Value = '''0.0270'
$ws.Range("E45").Value = '  +5.43%  '

$ws.Range("D46").Value = '''9.59'
$ws.Range("E46").Value = '  +4.38%  '

$ws.Range("D47").Value = '''16.13'
$ws.Range("E47").Value = '  -3.11%  '

$ws.Range("E48").Value = '  +0.04%  '

$ws.Range("D49").Value = '2.450.63'
$ws.Range("E49").Value = '  -3.24%  '

$ws.Range("D50").Value = '''70.86'
$ws.Range("E50").Value = '  +3.69%  '

$ws.Range("D51").Value = '''1.45'
$ws.Range("E51").Value = '  +9.97%  '
